# Update "想去人数" (want-to-go count) figures in the 江西-漫展信息 workbook
# to the values captured at the later site-generation run (456a3b4).
#
# Two worksheets hold (mostly) the same rows of exhibition data and both
# need the same numeric bumps applied to column F:
#   展览   (Worksheets("展览"))   -> rows 2,3,4,7,10,11,12,13,15
#   全部类型 (Worksheets("全部类型")) -> rows 2,3,4,7,11,12,13,14,16
# (全部类型 has one extra row - a performance entry - inserted at row 10,
#  which shifts the matching rows down by one from row 10 onward.)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4520
$wsExhibit.Range("F3").Value = 859
$wsExhibit.Range("F4").Value = 142
$wsExhibit.Range("F7").Value = 156
$wsExhibit.Range("F10").Value = 201
$wsExhibit.Range("F11").Value = 1392
$wsExhibit.Range("F12").Value = 30
$wsExhibit.Range("F13").Value = 2990
$wsExhibit.Range("F15").Value = 676

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4520
$wsAll.Range("F3").Value = 859
$wsAll.Range("F4").Value = 142
$wsAll.Range("F7").Value = 156
$wsAll.Range("F11").Value = 201
$wsAll.Range("F12").Value = 1392
$wsAll.Range("F13").Value = 30
$wsAll.Range("F14").Value = 2990
$wsAll.Range("F16").Value = 676
